$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tabla1")

# 1) Grow the table from A2:D6 to A2:F6 (adds ListColumns 5 and 6).
$lo.Resize($ws.Range("A2:F6"))

# 2) Move the existing "Examen" scores (column D) into the new
#    "Examen Practico" column (E) before we touch column D.
$ws.Range("D3:D6").Copy() | Out-Null
$ws.Range("E3:E6").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# 3) Fill in the brand new "Examen Teorico" scores in column F.
$ws.Range("F3").Value = 8
$ws.Range("F4").Value = 3.5
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 9

# 4) Give the new columns the same look (font/alignment) as the rest of
#    the data area by copying formatting from the neighboring column C.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D3:D6").PasteSpecial(-4122) | Out-Null
$ws.Range("E3:E6").PasteSpecial(-4122) | Out-Null
$ws.Range("F3:F6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 5) Re-title the headers. Column D becomes "Tarea 3" (no data entered
#    for it yet) while E/F pick up the exam names. Write E/F before D so
#    new shared-string entries land in the same order as the workbook
#    being reproduced.
$ws.Range("E2").Value = "Examen Practico"
$ws.Range("F2").Value = "Examen Teórico"
$ws.Range("D2").Value = "Tarea 3"

# Match the header cell formatting (it was already applied to D2, copy it
# across to the two brand new header cells).
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 6) Column D no longer holds scores - clear its body but keep the style.
$ws.Range("D3:D6").ClearContents() | Out-Null

# 7) Cosmetic touch-ups matching the authored workbook: widen the new/
#    edited columns and move the active selection.
$ws.Columns.Item(4).ColumnWidth = 28.67
$ws.Columns.Item(5).ColumnWidth = 31.83
$ws.Columns.Item(6).ColumnWidth = 34.67
$ws.Range("A4").Select() | Out-Null

$wb.Save()
